# Update (Analyze PO & Forecast)
# Updates the MyForecast column (D) on the "Forecast Comparison" sheet
# and the derived summary statistics (column B) on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# --- Forecast Comparison sheet: MyForecast values (column D, rows 2-17) ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$forecastValues = @{
    2  = 44
    3  = 46
    4  = 52
    5  = 55
    6  = 53
    7  = 53
    8  = 55
    9  = 50
    10 = 45
    11 = 53
    12 = 52
    13 = 55
    14 = 53
    15 = 52
    16 = 52
    17 = 51
}

foreach ($row in $forecastValues.Keys) {
    $wsForecast.Cells.Item($row, 4).Value = $forecastValues[$row]
}

# --- Summary sheet: updated forecast statistics (column B, rows 9-15) ---
# These cells store their values as TEXT (not numbers/dates), so force the
# cell to Text format before the write, then restore the default "Normal"
# style so no residual number-format is left behind on the cell.
$wsSummary = $wb.Worksheets.Item("Summary")

$summaryValues = @{
    9  = "821"
    10 = "408"
    11 = "197"
    12 = "55"
    13 = "2025-02-16"
    14 = "44"
    15 = "2025-01-26"
}

foreach ($row in $summaryValues.Keys) {
    $cell = $wsSummary.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryValues[$row]
    $cell.Style = "Normal"
}
